$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 73
$ws.Range("I2").Value = 207
$ws.Range("J2").Value = 783
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 221
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = 167
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 14
$ws.Range("S2").Value = 100
$ws.Range("T2").Value = 124
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 1219
$ws.Range("X2").Value = 1253
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 11
